$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.24"
$ws.Range("D4").Value = "'5.203"
$ws.Range("D5").Value = "'0.05606"
$ws.Range("D6").Value = "'3.370"
$ws.Range("D7").Value = "'6.395"
$ws.Range("D8").Value = "'0.8054"
$ws.Range("D9").Value = "'0.9740"
$ws.Range("D10").Value = "'0.1412"
$ws.Range("D11").Value = "'0.07295"
$ws.Range("D12").Value = "'0.03122"
$ws.Range("D13").Value = "'0.03058"
$ws.Range("D14").Value = "'0.09285"
$ws.Range("D15").Value = "'3.572"
$ws.Range("D16").Value = "'0.001651"
$ws.Range("D17").Value = "'0.04716"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006401"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.004986"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001042"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "UpBots"
$ws.Range("C22").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D22").Value = "'0.0003100"
$ws.Range("E22").Value = "21UpBotsUBXT"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.754"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.098"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Value = "'0.01120"
$ws.Range("E25").Value = "24OneONEBestin24h"
$ws.Range("D26").Value = "'0.3260"
$ws.Range("D27").Value = "'0.1271"
$ws.Range("D40").Value = "'0.03902"
$ws.Range("D41").Value = "'0.006876"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003400"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1037"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.008505"
$ws.Range("D45").Value = "'0.00005927"
$ws.Range("D47").Value = "'0.0005498"
$ws.Range("D48").Value = "'0.6826"
$ws.Range("D49").Value = "'0.09153"
$ws.Range("E49").Value = "48BOLOBOLO"
